$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump "Forandrad" (column C) date value by one day (46072 -> 46073) for all data rows (2-41)
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = 46073
}

# Refreshed listing data (columns A, B, F, G) for rows 5-41 (row 29 unaffected)

# Row 5
$ws.Cells.Item(5, 1).Value = "A 1379-2024"
$ws.Cells.Item(5, 2).Value = 45303.55193287037
$ws.Cells.Item(5, 6).Value = "Kommuner"
$ws.Cells.Item(5, 7).Value = 1.5

# Row 6
$ws.Cells.Item(6, 1).Value = "A 61627-2024"
$ws.Cells.Item(6, 2).Value = 45646.66263888889
$ws.Cells.Item(6, 6).Value = "Kommuner"
$ws.Cells.Item(6, 7).Value = 2.1

# Row 7
$ws.Cells.Item(7, 1).Value = "A 1486-2023"
$ws.Cells.Item(7, 2).Value = 44937
$ws.Cells.Item(7, 7).Value = 0.2

# Row 8
$ws.Cells.Item(8, 1).Value = "A 23301-2025"
$ws.Cells.Item(8, 2).Value = 45791.59498842592
$ws.Cells.Item(8, 7).Value = 2.9

# Row 9
$ws.Cells.Item(9, 1).Value = "A 23294-2025"
$ws.Cells.Item(9, 2).Value = 45791.58967592593
$ws.Cells.Item(9, 7).Value = 4.6

# Row 10
$ws.Cells.Item(10, 1).Value = "A 23295-2025"
$ws.Cells.Item(10, 2).Value = 45791.59071759259
$ws.Cells.Item(10, 7).Value = 2.2

# Row 11
$ws.Cells.Item(11, 1).Value = "A 42462-2025"
$ws.Cells.Item(11, 2).Value = 45905.45730324074
$ws.Cells.Item(11, 7).Value = 2.3

# Row 12
$ws.Cells.Item(12, 1).Value = "A 60718-2022"
$ws.Cells.Item(12, 2).Value = 44912.89109953704
$ws.Cells.Item(12, 7).Value = 1

# Row 13
$ws.Cells.Item(13, 1).Value = "A 45167-2025"
$ws.Cells.Item(13, 2).Value = 45919.49364583333
$ws.Cells.Item(13, 7).Value = 0.6

# Row 14
$ws.Cells.Item(14, 1).Value = "A 31486-2021"
$ws.Cells.Item(14, 2).Value = 44369.43783564815
$ws.Cells.Item(14, 7).Value = 4.8

# Row 15
$ws.Cells.Item(15, 1).Value = "A 45158-2025"
$ws.Cells.Item(15, 2).Value = 45919.48245370371
$ws.Cells.Item(15, 6).ClearContents()
$ws.Cells.Item(15, 7).Value = 2.7

# Row 16
$ws.Cells.Item(16, 1).Value = "A 45088-2025"
$ws.Cells.Item(16, 2).Value = 45919.37598379629
$ws.Cells.Item(16, 6).ClearContents()
$ws.Cells.Item(16, 7).Value = 0.9

# Row 17
$ws.Cells.Item(17, 1).Value = "A 57410-2022"
$ws.Cells.Item(17, 2).Value = 44896
$ws.Cells.Item(17, 7).Value = 7.5

# Row 18
$ws.Cells.Item(18, 1).Value = "A 43229-2024"
$ws.Cells.Item(18, 2).Value = 45567.88697916667
$ws.Cells.Item(18, 7).Value = 9.699999999999999

# Row 19
$ws.Cells.Item(19, 1).Value = "A 57000-2025"
$ws.Cells.Item(19, 2).Value = 45977
$ws.Cells.Item(19, 7).Value = 2.3

# Row 20
$ws.Cells.Item(20, 1).Value = "A 1103-2025"
$ws.Cells.Item(20, 2).Value = 45666
$ws.Cells.Item(20, 7).Value = 1.6

# Row 21
$ws.Cells.Item(21, 1).Value = "A 57655-2025"
$ws.Cells.Item(21, 2).Value = 45981.40369212963
$ws.Cells.Item(21, 7).Value = 3.4

# Row 22
$ws.Cells.Item(22, 1).Value = "A 769-2023"
$ws.Cells.Item(22, 2).Value = 44931
$ws.Cells.Item(22, 7).Value = 1.7

# Row 23
$ws.Cells.Item(23, 1).Value = "A 21379-2023"
$ws.Cells.Item(23, 2).Value = 45063.34819444444
$ws.Cells.Item(23, 7).Value = 5.8

# Row 24
$ws.Cells.Item(24, 1).Value = "A 34394-2025"
$ws.Cells.Item(24, 2).Value = 45846.58854166666
$ws.Cells.Item(24, 7).Value = 3.9

# Row 25
$ws.Cells.Item(25, 1).Value = "A 32577-2025"
$ws.Cells.Item(25, 2).Value = 45838
$ws.Cells.Item(25, 7).Value = 5.3

# Row 26
$ws.Cells.Item(26, 1).Value = "A 34508-2025"
$ws.Cells.Item(26, 2).Value = 45847.44315972222
$ws.Cells.Item(26, 7).Value = 2

# Row 27
$ws.Cells.Item(27, 1).Value = "A 2434-2026"
$ws.Cells.Item(27, 2).Value = 46036.86722222222
$ws.Cells.Item(27, 7).Value = 1.2

# Row 28
$ws.Cells.Item(28, 1).Value = "A 2433-2026"
$ws.Cells.Item(28, 2).Value = 46036.86631944445
$ws.Cells.Item(28, 7).Value = 2.6

# Row 30
$ws.Cells.Item(30, 1).Value = "A 34939-2025"
$ws.Cells.Item(30, 2).Value = 45849.58229166667
$ws.Cells.Item(30, 7).Value = 7.8

# Row 31
$ws.Cells.Item(31, 1).Value = "A 34963-2025"
$ws.Cells.Item(31, 2).Value = 45849.63219907408
$ws.Cells.Item(31, 7).Value = 1.1

# Row 32
$ws.Cells.Item(32, 1).Value = "A 59011-2025"
$ws.Cells.Item(32, 2).Value = 45987
$ws.Cells.Item(32, 7).Value = 2.7

# Row 33
$ws.Cells.Item(33, 1).Value = "A 7791-2023"
$ws.Cells.Item(33, 2).Value = 44973
$ws.Cells.Item(33, 7).Value = 3.1

# Row 34
$ws.Cells.Item(34, 1).Value = "A 34591-2022"
$ws.Cells.Item(34, 2).Value = 44795.3778587963
$ws.Cells.Item(34, 7).Value = 2.5

# Row 35
$ws.Cells.Item(35, 1).Value = "A 63664-2023"
$ws.Cells.Item(35, 2).Value = 45275.62074074074
$ws.Cells.Item(35, 7).Value = 3.2

# Row 36
$ws.Cells.Item(36, 1).Value = "A 60793-2023"
$ws.Cells.Item(36, 2).Value = 45260.6534375
$ws.Cells.Item(36, 7).Value = 0.7

# Row 37
$ws.Cells.Item(37, 1).Value = "A 57955-2024"
$ws.Cells.Item(37, 2).Value = 45631.56939814815
$ws.Cells.Item(37, 7).Value = 0.9

# Row 38
$ws.Cells.Item(38, 1).Value = "A 21141-2023"
$ws.Cells.Item(38, 2).Value = 45062
$ws.Cells.Item(38, 7).Value = 3.4

# Row 39
$ws.Cells.Item(39, 1).Value = "A 60392-2022"
$ws.Cells.Item(39, 2).Value = 44910
$ws.Cells.Item(39, 7).Value = 4.1

# Row 40
$ws.Cells.Item(40, 1).Value = "A 59192-2022"
$ws.Cells.Item(40, 2).Value = 44896
$ws.Cells.Item(40, 7).Value = 1.5

# Row 41
$ws.Cells.Item(41, 1).Value = "A 60717-2022"
$ws.Cells.Item(41, 2).Value = 44912.89078703704
$ws.Cells.Item(41, 7).Value = 0.9
